# 1 - CDP Network Audit _ Template.xlsx
# The "Seed Device 1:" / "Seed Device 2:" prompts on the Audit sheet are
# consolidated into a single "Seed Device:" prompt.
#   - A7 ("Seed Device 1:") -> "Seed Device:"
#   - A8 ("Seed Device 2:") -> cleared (kept blank, style preserved)
#   - Active selection on the Audit sheet moves from A17 to A8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Audit")

$ws.Range("A7").Value = "Seed Device:"
$ws.Range("A8").ClearContents()

[void]$ws.Range("A8").Select()
